# Add a new column (BH) of data to the "Prix Spot" sheet:
# header "12-aug" in BH1 (same style as the other date headers),
# and the corresponding numeric values in BH2:BH25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell - copy the style used by the existing date headers (column BG)
$ws.Range("BH1").Value = "12-aug"
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)  # xlPasteFormats

$values = @{
    2  = 93
    3  = 89.67
    4  = 85.41
    5  = 79.06
    6  = 69.67
    7  = 66.34
    8  = 75.52
    9  = 102.81
    10 = 98.69
    11 = 92.5
    12 = 66.34
    13 = 40.94
    14 = 8.380000000000001
    15 = 3
    16 = 5.15
    17 = 35.05
    18 = 63.8
    19 = 85.81999999999999
    20 = 96.78
    21 = 143.64
    22 = 154.93
    23 = 140.81
    24 = 118.41
    25 = 100.37
}

foreach ($row in $values.Keys) {
    $ws.Range("BH$row").Value = $values[$row]
}
